$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values: B1..E1 changed (16, 20, 16, 20)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON): B2 and C2 deleted, D2 and E2 updated
$ws.Range("B2").ClearContents() | Out-Null
$ws.Range("C2").ClearContents() | Out-Null
$ws.Range("D2").Value = 12.718423230545389
$ws.Range("E2").Value = 11.494860884208649

# Row 3 (STR): B3 deleted, C3 updated, D3 newly added, E3 updated
$ws.Range("B3").ClearContents() | Out-Null
$ws.Range("C3").Value = 9.5442078489037812
$ws.Range("D3").Value = 11.054580937996574
$ws.Range("E3").Value = 11.009267784812847

# Selection now covers only B1:E3 instead of the whole used range
$ws.Range("B1:E3").Select() | Out-Null
